# Set column L ("Diferencia Stock") equal to column K ("Stock Mínimo Objetivo")
# for each data row (3-44), and refresh the "Total_Ajuste_Stock" total in C58
# so it reflects the new sum of column L.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_7")

$total = 0
for ($row = 3; $row -le 44; $row++) {
    $kValue = $ws.Cells.Item($row, 11).Value2   # Column K (use Value2 - reliable getter)
    $ws.Cells.Item($row, 12).Value = $kValue    # Column L
    $total = $total + $kValue
}

# Update the "Total_Ajuste_Stock" summary cell to match the new column L total
$ws.Range("C58").Value = $total
